$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Amount:" values for AHUs and Pumps from 0 to 1
$ws.Range("J2").Value = 1
$ws.Range("M2").Value = 1

# Update the "Length:" value for Chillers from 1 to 3
$ws.Range("D3").Value = 3

# Move the selected/active cell to D4
$ws.Range("D4").Select()
